$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 4 new rows before old row 8 (HLR_00600/IT_00600 row).
#    This pushes the old rows 8,9,10 (HLR_00600/610/700) down to 12,13,14
#    and creates blank rows 8,9,10,11 for the new SPI/UART test cases.
# ------------------------------------------------------------------
$ws.Range("A8:A11").EntireRow.Insert()

# ------------------------------------------------------------------
# 2. Row heights (wrapped text needs taller rows for the longer text)
# ------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30

# ------------------------------------------------------------------
# 3. Cell content, in the same order the workbook author entered it
#    (this ordering controls how new shared-strings get appended).
# ------------------------------------------------------------------
$ws.Range("C5").Value = "Test if the ultrason sensor considered that it detects an obstacle under 8cm"
$ws.Range("C6").Value = "Test if the servomotor rotate in range of [+45°;-45°]"
$ws.Range("C7").Value = "Test if the infrared sensor considered that it detects a hole over 4cm"

$ws.Range("A8").Value = "HLR_00500"
$ws.Range("A9").Value = "HLR_00500"
$ws.Range("A10").Value = "HLR_00510"
$ws.Range("A11").Value = "HLR_00510"

$ws.Range("B8").Value = "IT_00500"
$ws.Range("B10").Value = "IT_00510"

$ws.Range("C8").Value = "Test if informations are sent by the user to MSP430G2553 with UART"
$ws.Range("C9").Value = "Test if informations are sent by the MSP430G2553 to user with UART"
$ws.Range("C10").Value = "Test if informations are sent by the MSP430G2553 to MSP430G2231 with SPI"
$ws.Range("C11").Value = "Test if informations are sent by the MSP430G2231 to MSP430G2553 with SPI"

$ws.Range("B9").Value = "IT_00501"
$ws.Range("B11").Value = "IT_00511"

$ws.Range("C12").Value = "Test if the bot start when the command is sent"
$ws.Range("C13").Value = "Test if the bot stop when the command is sent"
$ws.Range("C14").Value = "Test if sensors data are display when user send command"

$ws.Range("E5").Value = "Ultrason sensor returns 1 (obstacle detected)"
$ws.Range("E7").Value = "Infrared sensor returns 1 (hole detected)"
$ws.Range("E6").Value = "Servomotor rotate in range of [+45°;-45°]"
$ws.Range("E9").Value = "User receive informations sent by the MSP430G2553"
$ws.Range("E8").Value = "MSP430G2553 receive informations sent but the user"
$ws.Range("E10").Value = "MSP430G2231 receive the informations sent by the MSP430G2553"
$ws.Range("E11").Value = "MSP430G2553 receive the informations sent by the MSP430G2231"
$ws.Range("E12").Value = "The bot start"
$ws.Range("E13").Value = "The bot stop"
$ws.Range("E14").Value = "Sensors data are display"

# ------------------------------------------------------------------
# 4. Grow the table to cover the new rows
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F14"))

# ------------------------------------------------------------------
# 5. Column widths: column E becomes wider than the rest
# ------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 37.2857

# ------------------------------------------------------------------
# 6. Extra blank formatted rows at the bottom of the sheet
# ------------------------------------------------------------------
$tailRng = $ws.Range("A39:F42")
$tailRng.HorizontalAlignment = -4108
$tailRng.VerticalAlignment = -4108
$tailRng.WrapText = $true

# ------------------------------------------------------------------
# 7. Conditional formatting on the Pass/Fail column
# ------------------------------------------------------------------
$cfRng = $ws.Range("F2:F14")
$fcTemp = $cfRng.FormatConditions.Add(9, 0, $null, $null, "PASS")
$fcTemp.Interior.Color = 5287936
$fcTemp.Delete()

$fcFail = $cfRng.FormatConditions.Add(9, 0, $null, $null, "FAIL")
$fcFail.Interior.Color = 255

$fcPass = $cfRng.FormatConditions.Add(9, 0, $null, $null, "PASS")
$fcPass.Interior.Color = 5287936

# ------------------------------------------------------------------
# 8. Selection / view state
# ------------------------------------------------------------------
$ws.Range("F4").Select()
